# Update liquefaction probability formula (column P) to drop the *100 scaling
# and the now-unneeded B>1700 branch, lowering the floor from 0.001 to 0.00001
# so that computed probabilities are returned in decimal/fraction form instead
# of percent, and to avoid the prior formula returning values > 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hazus2020_with_ZhuEtal2017")

# Row 2 has its own (non-shared) formula.
$ws.Range("P2").Formula = "=IF(I2<0.1,0.00001,IF(A2>620,0.00001,MAX(N2/L2/M2*O2,0.00001)))"

# Rows 3:25 share one formula (relative references adjust per row automatically).
$ws.Range("P3:P25").Formula = "=IF(I3<0.1,0.00001,IF(A3>620,0.00001,MAX(N3/L3/M3*O3,0.00001)))"

# Recalculate so the dependent Q column (and any other dependents) reflect the
# new P values.
$excel.CalculateFullRebuild()

# Match the author's final selection state (P2:P25 instead of Q2:Q25).
$ws.Range("P2:P25").Select()
